# GateInward with Balloon on SaveEvent
#
# Adds a new task row (row 14) describing "Change from ASCII to english
# (For Item Code Des)" task details, adds a new "Ord_Hist, Ord_Det" table
# reference to row 13 (B13 / I13), and updates the sheet selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 13: add the "Ord_Hist, Ord_Det" tables note (column B) and
#     duplicate the "Done" marker into the new Item Code Description
#     column (I). Copy formatting from sibling cells so styles match. ---

$ws.Range("A13").Copy()
$ws.Range("B13").PasteSpecial(-4122, $false, $false, $false)
$ws.Range("B13").Value = "Ord_Hist, Ord_Det"

$ws.Range("H13").Copy()
$ws.Range("I13").PasteSpecial(-4122, $false, $false, $false)
$ws.Range("I13").Value = "Done"

# --- Row 14: fill in the rest of the "Change from ASCII to english" task:
#     Start Date, Status and the Done markers across F:J. ---

$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122, $false, $false, $false)
$ws.Range("C14").Value = 41656

$ws.Range("F13").Copy()
$ws.Range("F14").PasteSpecial(-4122, $false, $false, $false)
$ws.Range("F14").Value = "Completed & Revised"

$ws.Range("G13").Copy()
$ws.Range("G14").PasteSpecial(-4122, $false, $false, $false)
$ws.Range("G14").Value = "Done"

$ws.Range("H13").Copy()
$ws.Range("H14").PasteSpecial(-4122, $false, $false, $false)
$ws.Range("H14").Value = "Done"

$ws.Range("I13").Copy()
$ws.Range("I14").PasteSpecial(-4122, $false, $false, $false)
$ws.Range("I14").Value = "Done"

$ws.Range("J13").Copy()
$ws.Range("J14").PasteSpecial(-4122, $false, $false, $false)
$ws.Range("J14").Value = "Done"

# --- Sheet view: move the selection down onto the newly edited row. ---

$ws.Range("K13").Select()
